# [Ajout ICONES + creation page RECLAMATION]
# Resize/reposition the "Rectangle 4" shape (plateFORME logo, shadowed variant)
# and bump up its run font size from 44pt to 80pt (endParaRPr follow-up size
# bump from 20pt to 44pt is a trailing/invisible paragraph-mark property not
# reachable through the PowerPoint automation surface).

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(1)
$sh = $s.Shapes.Item(2)   # "Rectangle 4"

# New position/size (EMU target: off 4131733,4157133  ext 7984067,1159934).
# Point values below are chosen so the host's pt -> EMU conversion lands
# exactly on the target EMU amounts.
$sh.Left   = 325.3332977866142
$sh.Top    = 327.3332977866142
$sh.Width  = 628.6667176133859
$sh.Height = 91.33338552677165

# Bump the visible run's font size from 44pt to 80pt.
$tf = $sh.TextFrame
$tr = $tf.TextRange
$tr.Font.Size = 80
